# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Overview sheet: status columns widen (auto-fit for longer status text)
#  - zh-cn / de-de sheets: "Latest Target File" / "Latest Handback File" /
#    "Latest Handback DateTime" columns get populated for the rows that were
#    just handed back, and their hyperlinks + column widths follow suit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Shared-string-backed text updates (Status column + datetime values)
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")
$overview = $wb.Worksheets.Item("Overview")

# "Ready for handoff" -> "Handed back: in sync with en-US" everywhere it is used
# (Overview!E2:F3 and the Status column (C) on the zh-cn / de-de sheets).
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"

# zh-cn "Latest Handback DateTime" (K2/K3) moves from the epoch placeholder
# to a real timestamp.
$zhcn.Range("K2").Value = "2016-09-03 18:49:39"
$zhcn.Range("K3").Value = "2016-09-03 18:49:39"

# ---------------------------------------------------------------------
# 2. zh-cn sheet: populate "Latest Target File" (I) / "Latest Handback
#    File" (J) for rows 2 and 3, with I2/I3 becoming hyperlinks to a.md.
# ---------------------------------------------------------------------
$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/233156f787e58676f4e2988a5ee13994fee42dd8/e2e/a.md", "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/233156f787e58676f4e2988a5ee13994fee42dd8/e2e/a.md", "", "", "a.md")
$zhcn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/233156f787e58676f4e2988a5ee13994fee42dd8/e2e/b.md", "", "", "b.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/233156f787e58676f4e2988a5ee13994fee42dd8/e2e/a.md", "", "", "a.md")
$zhcn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

# ---------------------------------------------------------------------
# 3. de-de sheet: same shape of update, plus "Latest Handback DateTime"
#    (K) gets the freshly generated handback timestamp.
# ---------------------------------------------------------------------
$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/233156f787e58676f4e2988a5ee13994fee42dd8/e2e/a.md", "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/233156f787e58676f4e2988a5ee13994fee42dd8/e2e/a.md", "", "", "a.md")
$dede.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K2").Value = "2016-09-03 18:49:46"
$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/233156f787e58676f4e2988a5ee13994fee42dd8/e2e/b.md", "", "", "b.md")
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/233156f787e58676f4e2988a5ee13994fee42dd8/e2e/a.md", "", "", "a.md")
$dede.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K3").Value = "2016-09-03 18:49:46"

# ---------------------------------------------------------------------
# 4. Column widths widen to fit the new, longer content. The engine
#    quantizes ColumnWidth to a coarse grid, so we pick the input value
#    that lands closest to the canonical widths used by the source file.
# ---------------------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 29.1666666666667
$overview.Columns.Item(6).ColumnWidth = 29.1666666666667

$zhcn.Columns.Item(3).ColumnWidth = 29.1666666666667
$zhcn.Columns.Item(10).ColumnWidth = 39.1666666666667

$dede.Columns.Item(3).ColumnWidth = 29.1666666666667
$dede.Columns.Item(10).ColumnWidth = 39.1666666666667

Write-Output "Handback report generated"
